$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 - copy H1's formatting (bold/bordered/centered header
# style) so the new header cells match the rest of row 1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 (start inning) / IF (finish inning) values per game row, derived from IP
# (column H): IF = I0 + IP - 1
$data = @{
    2  = @(1, 5)
    3  = @(1, 3)
    4  = @(1, 1)
    5  = @(1, 4)
    6  = @(1, 2)
    7  = @(1, 3)
    8  = @(1, 5)
    9  = @(1, 3)
    10 = @(1, 6)
    11 = @(1, 7)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 5)
    15 = @(1, 4)
    16 = @(1, 3)
    17 = @(1, 6)
    18 = @(7, 8)
    19 = @(1, 6)
    20 = @(1, 6)
    21 = @(1, 6)
    22 = @(1, 5)
    23 = @(1, 4)
    24 = @(4, 6)
    25 = @(8, 8)
    26 = @(6, 7)
    27 = @(4, 5)
    28 = @(1, 2)
    29 = @(3, 4)
}

foreach ($row in $data.Keys) {
    $pair = $data[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
